# Updated cryptos list - apply new price/volume data per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "61.631.92"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +3.67%  "

# Row 3
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.075.64"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +2.43%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "575.68"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +1.83%  "

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "141.60"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +1.41%  "

# Row 7
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "3.064.22"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +2.33%  "

# Row 9
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.525"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.74%  "

# Row 10
$ws.Range("E10").Value = "  +3.36%  "

# Row 11
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "5.48"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +8.88%  "

# Row 12
$ws.Range("E12").Value = "  +1.02%  "

# Row 13
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.0000239"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +3.32%  "

# Row 14
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "34.96"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +2.99%  "

# Row 15
$ws.Range("E15").Value = "  +0.05%  "

# Row 16
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "3.582.02"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +2.41%  "

# Row 17
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "7.24"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +3.18%  "

# Row 18
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "3.068.50"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +2.30%  "

# Row 19
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "61.565.84"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +3.78%  "

# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "448.78"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +4.23%  "

# Row 21
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "13.92"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +2.31%  "

# Row 22
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.730"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +2.36%  "

# Row 23
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "7.42"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +3.13%  "

# Row 24
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "13.52"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.22%  "

# Row 25
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "81.89"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +1.39%  "

# Row 26
$ws.Range("E26").Value = "  +0.19%  "

# Row 27
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "2.23"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +4.22%  "

# Row 28
$ws.Range("B28").Value = "FirstDigitalUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +0.11%  "

# Row 29
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.63"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +2.94%  "

# Row 30
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "8.02"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +2.59%  "

# Row 31
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "6.68"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +8.04%  "

# Row 32
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "26.54"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +3.07%  "

# Row 33
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.109"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +9.30%  "

# Row 34
$ws.Range("E34").Value = "  +1.76%  "

# Row 35
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0792"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +2.31%  "

# Row 36
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "6.05"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +4.60%  "

# Row 37
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.16"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +3.40%  "

# Row 38
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "50.18"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +1.29%  "

# Row 39
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "2.97"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +7.18%  "

# Row 40
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "8.81"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +1.28%  "

# Row 41
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "422.65"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +3.69%  "

# Row 42
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.0368"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +4.14%  "

# Row 43
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "2.771.19"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +0.34%  "

# Row 44
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.268"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +7.06%  "

# Row 45
$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.108"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -0.41%  "

# Row 46
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "35.65"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +8.51%  "

# Row 47
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "2.09"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +3.04%  "

# Row 48
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "124.91"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +0.26%  "

# Row 49
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +0.01%  "

# Row 50
$ws.Range("E50").Value = "  +0.98%  "

# Row 51
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "23.82"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +0.80%  "
